# "Form Sorting, except titles completed. Added icons"
#
# The row for Sno=8 ("dffc" choice entry) is removed from the data table;
# every row below it shifts up by one, the used range shrinks from
# A1:F11 to A1:F10, the AutoFilter (which was applied only to column D)
# is turned off, its hidden _FilterDatabase defined name is updated to
# the new (smaller) range, and the active selection moves to D17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row that held Sno = 8 ("dffc"). Excel shifts the
# remaining rows up and compacts the shared-string table automatically.
$ws.Rows(9).Delete() | Out-Null

# Data no longer needs the AutoFilter - turn it off (removes the
# <autoFilter> element from the sheet).
$ws.AutoFilterMode = $false

# The hidden _xlnm._FilterDatabase defined name survives the filter
# removal but still has to be re-pointed at the smaller D1:D10 range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$D`$1:`$D`$10"
    }
}

# Leave the cursor where the user left it when they saved.
$ws.Range("D17").Select() | Out-Null
